# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across ALC, ARM, CRP, CUL, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1803.4193
$ws.Range("I40").Value = 1650.5
$ws.Range("J40").Value = 2081.4546
$ws.Range("K40").Value = 1650.5
$ws.Range("L40").Value = 2081.4546
$ws.Range("M40").Value = -1475.5
$ws.Range("N40").Value = -2431.4546

# Row 62
$ws.Range("H62").Value = 4804.7617
$ws.Range("I62").Value = 2876.6667
$ws.Range("J62").Value = 6250.8335
$ws.Range("K62").Value = 2876.6667
$ws.Range("L62").Value = 6250.8335
$ws.Range("M62").Value = -2252.6667
$ws.Range("N62").Value = -7498.8335

# Row 64
$ws.Range("H64").Value = 3405.3914
$ws.Range("I64").Value = 3131.9167
$ws.Range("J64").Value = 3703.7273
$ws.Range("K64").Value = 3131.9167
$ws.Range("L64").Value = 3703.7273
$ws.Range("M64").Value = -2883.9167
$ws.Range("N64").Value = -4199.7273

# Row 65
$ws.Range("H65").Value = 4804.7617
$ws.Range("I65").Value = 2876.6667
$ws.Range("J65").Value = 6250.8335
$ws.Range("K65").Value = 14383.3335
$ws.Range("L65").Value = 31254.1675
$ws.Range("M65").Value = -11263.3335
$ws.Range("N65").Value = -37494.1675

# Row 67
$ws.Range("H67").Value = 3405.3914
$ws.Range("I67").Value = 3131.9167
$ws.Range("J67").Value = 3703.7273
$ws.Range("K67").Value = 3131.9167
$ws.Range("L67").Value = 3703.7273
$ws.Range("M67").Value = -2273.9167
$ws.Range("N67").Value = -5419.7273

# Row 82
$ws.Range("H82").Value = 1012.6667
$ws.Range("I82").Value = 1012.6667
$ws.Range("K82").Value = 3038.0001
$ws.Range("M82").Value = -2632.0001

# Row 85
$ws.Range("H85").Value = 1012.6667
$ws.Range("I85").Value = 1012.6667
$ws.Range("K85").Value = 3038.0001
$ws.Range("M85").Value = -1634.0001

# Row 98
$ws.Range("H98").Value = 3999
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 3999
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 3999
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -6995

# Row 122
$ws.Range("H122").Value = 3999
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11997
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16897

# Row 138
$ws.Range("H138").Value = 3450283.8
$ws.Range("I138").Value = 1146.1136
$ws.Range("J138").Value = 14290431
$ws.Range("K138").Value = 3438.3408
$ws.Range("L138").Value = 42871293
$ws.Range("M138").Value = 1701.6592
$ws.Range("N138").Value = -42881573

# Row 139
$ws.Range("H139").Value = 19999
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Row 141
$ws.Range("H141").Value = 2375
$ws.Range("I141").Value = 2375
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7125
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1945
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 9196.77
$ws.Range("I21").Value = 4365.4287
$ws.Range("J21").Value = 14833.333
$ws.Range("K21").Value = 4365.4287
$ws.Range("L21").Value = 14833.333
$ws.Range("M21").Value = -3991.4287
$ws.Range("N21").Value = -15581.333

# Row 74
$ws.Range("H74").Value = 39841.652
$ws.Range("I74").Value = 84404.586
$ws.Range("J74").Value = 1644.8572
$ws.Range("K74").Value = 84404.586
$ws.Range("L74").Value = 1644.8572
$ws.Range("M74").Value = -83530.586
$ws.Range("N74").Value = -3392.8572

# Row 77
$ws.Range("H77").Value = 39841.652
$ws.Range("I77").Value = 84404.586
$ws.Range("J77").Value = 1644.8572
$ws.Range("K77").Value = 422022.93
$ws.Range("L77").Value = 8224.286
$ws.Range("M77").Value = -417654.93
$ws.Range("N77").Value = -16960.286

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5557.625
$ws.Range("I58").Value = 7139.706
$ws.Range("J58").Value = 1715.4286
$ws.Range("K58").Value = 7139.706
$ws.Range("L58").Value = 1715.4286
$ws.Range("M58").Value = -6936.706
$ws.Range("N58").Value = -2121.4286

# Row 70
$ws.Range("H70").Value = 33333.332
$ws.Range("J70").Value = 33333.332
$ws.Range("L70").Value = 33333.332
$ws.Range("N70").Value = -33963.332

# Row 73
$ws.Range("H73").Value = 33333.332
$ws.Range("J73").Value = 33333.332
$ws.Range("L73").Value = 33333.332
$ws.Range("N73").Value = -35517.332

# Row 136
$ws.Range("H136").Value = 5557.625
$ws.Range("I136").Value = 7139.706
$ws.Range("J136").Value = 1715.4286
$ws.Range("K136").Value = 21419.118
$ws.Range("L136").Value = 5146.2858
$ws.Range("M136").Value = -18869.118
$ws.Range("N136").Value = -10246.2858

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 4200.125
$ws.Range("J80").Value = 4200.125
$ws.Range("L80").Value = 12600.375
$ws.Range("N80").Value = -14472.375

# Row 83
$ws.Range("H83").Value = 4200.125
$ws.Range("J83").Value = 4200.125
$ws.Range("L83").Value = 37801.125
$ws.Range("N83").Value = -47161.125

# Row 107
$ws.Range("H107").Value = 2281001.8
$ws.Range("I107").Value = 261.1111
$ws.Range("J107").Value = 4333668.5
$ws.Range("K107").Value = 783.3333
$ws.Range("L107").Value = 13001005.5
$ws.Range("M107").Value = 1136.6667
$ws.Range("N107").Value = -13004845.5

# Row 131
$ws.Range("H131").Value = 850.5700000000001
$ws.Range("J131").Value = 893.1889
$ws.Range("L131").Value = 2679.5667
$ws.Range("N131").Value = -12759.5667

$ws = $wb.Worksheets.Item("LTW")
# Row 57
$ws.Range("H57").Value = 16515.334
$ws.Range("I57").Value = 9500
$ws.Range("K57").Value = 9500
$ws.Range("M57").Value = -8934

# Row 132
$ws.Range("H132").Value = 3633.5938
$ws.Range("I132").Value = 3108.348
$ws.Range("J132").Value = 4975.8887
$ws.Range("K132").Value = 9325.044
$ws.Range("L132").Value = 14927.6661
$ws.Range("M132").Value = -6795.044
$ws.Range("N132").Value = -19987.6661

# Row 136
$ws.Range("H136").Value = 17545444
$ws.Range("I136").Value = 19609350
$ws.Range("K136").Value = 58828050
$ws.Range("M136").Value = -58825500

$ws = $wb.Worksheets.Item("WVR")
# Row 55
$ws.Range("H55").Value = 7791.75
$ws.Range("I55").Value = 3274.6667
$ws.Range("J55").Value = 12308.833
$ws.Range("K55").Value = 3274.6667
$ws.Range("L55").Value = 12308.833
$ws.Range("M55").Value = -2997.6667
$ws.Range("N55").Value = -12862.833

# Row 132
$ws.Range("H132").Value = 6455516
$ws.Range("I132").Value = 8337271
$ws.Range("J132").Value = 3786
$ws.Range("K132").Value = 25011813
$ws.Range("L132").Value = 11358
$ws.Range("M132").Value = -25009283
$ws.Range("N132").Value = -16418

# Row 136
$ws.Range("H136").Value = 62504972
$ws.Range("I136").Value = 71433680
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 214301040
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -214298490
$ws.Range("N136").Value = -17100
